$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "mg / (1l)"
$ws.Range("F2").Value = "mg / (1l)"
$ws.Range("G2").Value = "mmol / (1l)"
$ws.Range("H2").Value = "mg / (1l)"
$ws.Range("I2").Value = "mg / (1l)"
$ws.Range("J2").Value = "mg / (1l)"
$ws.Range("K2").Value = "mg / (1l)"
$ws.Range("L2").Value = "mg / (1l)"
$ws.Range("O2").Value = "ug / (1l)"
$ws.Range("E4").Value = "<1"
$ws.Range("F4").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("L4").Value = "<1"
$ws.Range("O4").Value = 1000
$ws.Range("E5").Value = "<1"
$ws.Range("F5").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("L5").Value = -1
$ws.Range("O5").Value = 1000
$ws.Range("E6").Value = "<1"
$ws.Range("F6").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("O6").Value = 1000
$ws.Range("E7").Value = "<1"
$ws.Range("F7").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 1
$ws.Range("L7").Value = 1
$ws.Range("O7").Value = 1000
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 1
$ws.Range("L8").Value = 1
$ws.Range("O8").Value = 1000
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 1
$ws.Range("L9").Value = "<1"
$ws.Range("O9").Value = 1000
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 1
$ws.Range("L10").Value = -1
$ws.Range("O10").Value = 1000
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 1
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 1
$ws.Range("O11").Value = 1000
$ws.Range("E12").Value = 0.3261379151327745
$ws.Range("F12").Value = 225.8966630056657
$ws.Range("G12").Value = 21.73653150166828
$ws.Range("H12").Value = 0.7764908977817396
$ws.Range("I12").Value = 0.4674349206032192
$ws.Range("J12").Value = 0.3337927559737088
$ws.Range("L12").Value = 1000
$ws.Range("O12").Value = 1000
$ws.Range("E13").Value = 0.3261379151327745
$ws.Range("F13").Value = 225.8966630056657
$ws.Range("G13").Value = 21.73653150166828
$ws.Range("H13").Value = 0.7764908977817396
$ws.Range("I13").Value = 0.4674349206032192
$ws.Range("J13").Value = 0.3337927559737088
$ws.Range("L13").Value = 1000
$ws.Range("O13").Value = 1000
$ws.Range("E14").Value = "<0.32613791"
$ws.Range("F14").Value = 225.8966630056657
$ws.Range("G14").Value = 21.73653150166828
$ws.Range("H14").Value = 0.7764908977817396
$ws.Range("I14").Value = 0.4674349206032192
$ws.Range("J14").Value = 0.3337927559737088
$ws.Range("L14").Value = 1000
$ws.Range("O14").Value = 1000
$ws.Range("E15").Value = "<0.32613791"
$ws.Range("F15").Value = 225.8966630056657
$ws.Range("G15").Value = 21.73653150166828
$ws.Range("H15").Value = 0.7764908977817396
$ws.Range("I15").Value = 0.4674349206032192
$ws.Range("J15").Value = 0.3337927559737088
$ws.Range("L15").Value = 1000
$ws.Range("O15").Value = 1000
$ws.Range("E16").Value = "<0.32613791"
$ws.Range("F16").Value = 225.8966630056657
$ws.Range("G16").Value = 21.73653150166828
$ws.Range("H16").Value = 0.7764908977817396
$ws.Range("I16").Value = 0.4674349206032192
$ws.Range("J16").Value = 0.3337927559737088
$ws.Range("L16").Value = 1000
$ws.Range("O16").Value = 1000
